$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Invert the binary label values in column B (rows 2-93): 0 -> 1, 1 -> 0
for ($r = 2; $r -le 93; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cur = $cell.Value()
    if ($cur -eq 0) {
        $cell.Value = 1
    } else {
        $cell.Value = 0
    }
}

# Reflect the recorded selection state: the full data range is selected
[void]$ws.Range("A1:B93").Select()
